# Auto-generated Excel COM-interop script
# Updates market-data derived cells (profit calc columns H-N) across all 8 sheets
# to reflect a scheduled-runner data refresh, per commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC: set 120 cell value(s)
$ws.Range("H17").Value = 744.6667
$ws.Range("J17").Value = 715.63635
$ws.Range("L17").Value = 2146.90905
$ws.Range("N17").Value = -2482.90905
$ws.Range("H33").Value = 590.875
$ws.Range("I33").Value = 303.81818
$ws.Range("K33").Value = 303.81818
$ws.Range("M33").Value = -74.81817999999998
$ws.Range("H58").Value = 451.26666
$ws.Range("I58").Value = 269.2143
$ws.Range("K58").Value = 807.6428999999999
$ws.Range("M58").Value = -657.6428999999999
$ws.Range("H62").Value = 19882.412
$ws.Range("I62").Value = 17364.273
$ws.Range("K62").Value = 17364.273
$ws.Range("M62").Value = -16740.273
$ws.Range("H64").Value = 6671.409
$ws.Range("I64").Value = 3732.3333
$ws.Range("K64").Value = 3732.3333
$ws.Range("M64").Value = -3484.3333
$ws.Range("H65").Value = 19882.412
$ws.Range("I65").Value = 17364.273
$ws.Range("K65").Value = 86821.36500000001
$ws.Range("M65").Value = -83701.36500000001
$ws.Range("H67").Value = 6671.409
$ws.Range("I67").Value = 3732.3333
$ws.Range("K67").Value = 3732.3333
$ws.Range("M67").Value = -2874.3333
$ws.Range("H69").Value = 8551.5
$ws.Range("I69").Value = 5099.75
$ws.Range("K69").Value = 15299.25
$ws.Range("M69").Value = -14425.25
$ws.Range("H72").Value = 8551.5
$ws.Range("I72").Value = 5099.75
$ws.Range("K72").Value = 45897.75
$ws.Range("M72").Value = -41529.75
$ws.Range("H74").Value = 7301.9165
$ws.Range("I74").Value = 3454.75
$ws.Range("K74").Value = 3454.75
$ws.Range("M74").Value = -2518.75
$ws.Range("H77").Value = 7301.9165
$ws.Range("I77").Value = 3454.75
$ws.Range("K77").Value = 17273.75
$ws.Range("M77").Value = -12593.75
$ws.Range("H80").Value = 6945600
$ws.Range("I80").Value = 12346335
$ws.Range("K80").Value = 37039005
$ws.Range("M80").Value = -37038007
$ws.Range("H83").Value = 6945600
$ws.Range("I83").Value = 12346335
$ws.Range("K83").Value = 111117015
$ws.Range("M83").Value = -111112023
$ws.Range("H86").Value = 2297.0715
$ws.Range("I86").Value = 1487.25
$ws.Range("J86").Value = 2621
$ws.Range("K86").Value = 1487.25
$ws.Range("L86").Value = 2621
$ws.Range("M86").Value = -364.25
$ws.Range("N86").Value = -4867
$ws.Range("H89").Value = 2297.0715
$ws.Range("I89").Value = 1487.25
$ws.Range("J89").Value = 2621
$ws.Range("K89").Value = 7436.25
$ws.Range("L89").Value = 13105
$ws.Range("M89").Value = -1820.25
$ws.Range("N89").Value = -24337
$ws.Range("H98").Value = 1912.8182
$ws.Range("I98").Value = 1961.6666
$ws.Range("J98").Value = 1693
$ws.Range("K98").Value = 1961.6666
$ws.Range("L98").Value = 1693
$ws.Range("M98").Value = -463.6666
$ws.Range("N98").Value = -4689
$ws.Range("H101").Value = 1219.3572
$ws.Range("I101").Value = 935.8182
$ws.Range("J101").Value = 2259
$ws.Range("K101").Value = 2807.4546
$ws.Range("L101").Value = 6777
$ws.Range("M101").Value = -1185.4546
$ws.Range("N101").Value = -10021
$ws.Range("H112").Value = 7460.2954
$ws.Range("J112").Value = 7460.2954
$ws.Range("L112").Value = 22380.8862
$ws.Range("N112").Value = -24596.8862
$ws.Range("H122").Value = 1912.8182
$ws.Range("I122").Value = 1961.6666
$ws.Range("J122").Value = 1693
$ws.Range("K122").Value = 5884.9998
$ws.Range("L122").Value = 5079
$ws.Range("M122").Value = -3434.9998
$ws.Range("N122").Value = -9979
$ws.Range("H125").Value = 2334.6
$ws.Range("J125").Value = 2410.25
$ws.Range("L125").Value = 21692.25
$ws.Range("N125").Value = -26612.25
$ws.Range("H127").Value = 2982.8572
$ws.Range("I127").Value = 2982.8572
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 8948.571599999999
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -3988.571599999999
$ws.Range("H132").Value = 68902.484
$ws.Range("I132").Value = 72320.664
$ws.Range("K132").Value = 216961.992
$ws.Range("M132").Value = -214431.992
$ws.Range("H137").Value = 1356562.6
$ws.Range("I137").Value = 4005.4
$ws.Range("K137").Value = 12016.2
$ws.Range("M137").Value = -9466.200000000001
$ws.Range("H138").Value = 1771.081
$ws.Range("I138").Value = 1084.28
$ws.Range("J138").Value = 3201.9167
$ws.Range("K138").Value = 3252.84
$ws.Range("L138").Value = 9605.750100000001
$ws.Range("M138").Value = 1887.16
$ws.Range("N138").Value = -19885.7501
$ws.Range("H141").Value = 1317.7368
$ws.Range("I141").Value = 1168.7222
$ws.Range("K141").Value = 3506.1666
$ws.Range("M141").Value = 1673.8334
# ALC: clear 1 cell(s)
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM: set 65 cell value(s)
$ws.Range("H32").Value = 6669949
$ws.Range("I32").Value = 6669949
$ws.Range("K32").Value = 6669949
$ws.Range("M32").Value = -6669662
$ws.Range("H61").Value = 8334682.5
$ws.Range("I61").Value = 8334682.5
$ws.Range("K61").Value = 8334682.5
$ws.Range("M61").Value = -8334470.5
$ws.Range("H63").Value = 37300.332
$ws.Range("I63").Value = 1894
$ws.Range("J63").Value = 55003.5
$ws.Range("K63").Value = 1894
$ws.Range("L63").Value = 55003.5
$ws.Range("M63").Value = -1208
$ws.Range("N63").Value = -56375.5
$ws.Range("H66").Value = 37300.332
$ws.Range("I66").Value = 1894
$ws.Range("J66").Value = 55003.5
$ws.Range("K66").Value = 9470
$ws.Range("L66").Value = 275017.5
$ws.Range("M66").Value = -6038
$ws.Range("N66").Value = -281881.5
$ws.Range("H74").Value = 5757.0713
$ws.Range("I74").Value = 5649.5
$ws.Range("J74").Value = 5775
$ws.Range("K74").Value = 5649.5
$ws.Range("L74").Value = 5775
$ws.Range("M74").Value = -4775.5
$ws.Range("N74").Value = -7523
$ws.Range("H77").Value = 5757.0713
$ws.Range("I77").Value = 5649.5
$ws.Range("J77").Value = 5775
$ws.Range("K77").Value = 28247.5
$ws.Range("L77").Value = 28875
$ws.Range("M77").Value = -23879.5
$ws.Range("N77").Value = -37611
$ws.Range("H102").Value = 32167.916
$ws.Range("I102").Value = 42225.332
$ws.Range("J102").Value = 1995.6666
$ws.Range("K102").Value = 42225.332
$ws.Range("L102").Value = 1995.6666
$ws.Range("M102").Value = -40603.332
$ws.Range("N102").Value = -5239.6666
$ws.Range("H122").Value = 1544.7778
$ws.Range("I122").Value = 1401.5
$ws.Range("J122").Value = 1659.4
$ws.Range("K122").Value = 4204.5
$ws.Range("L122").Value = 4978.200000000001
$ws.Range("M122").Value = -1754.5
$ws.Range("N122").Value = -9878.200000000001
$ws.Range("H132").Value = 1484638.4
$ws.Range("I132").Value = 1978718.1
$ws.Range("J132").Value = 2399
$ws.Range("K132").Value = 5936154.300000001
$ws.Range("L132").Value = 7197
$ws.Range("M132").Value = -5933624.300000001
$ws.Range("N132").Value = -12257
$ws.Range("H133").Value = 98000
$ws.Range("J133").Value = 98000
$ws.Range("L133").Value = 98000
$ws.Range("N133").Value = -103060
$ws.Range("H136").Value = 8334682.5
$ws.Range("I136").Value = 8334682.5
$ws.Range("K136").Value = 25004047.5
$ws.Range("M136").Value = -25001497.5

$ws = $wb.Worksheets.Item("BSM")
# BSM: set 33 cell value(s)
$ws.Range("H82").Value = 57022.816
$ws.Range("J82").Value = 73104.625
$ws.Range("L82").Value = 73104.625
$ws.Range("N82").Value = -73870.625
$ws.Range("H85").Value = 57022.816
$ws.Range("J85").Value = 73104.625
$ws.Range("L85").Value = 73104.625
$ws.Range("N85").Value = -75756.625
$ws.Range("H86").Value = 5199
$ws.Range("I86").Value = 5329.6665
$ws.Range("J86").Value = 5003
$ws.Range("K86").Value = 5329.6665
$ws.Range("L86").Value = 5003
$ws.Range("M86").Value = -4206.6665
$ws.Range("N86").Value = -7249
$ws.Range("H89").Value = 5199
$ws.Range("I89").Value = 5329.6665
$ws.Range("J89").Value = 5003
$ws.Range("K89").Value = 26648.3325
$ws.Range("L89").Value = 25015
$ws.Range("M89").Value = -21032.3325
$ws.Range("N89").Value = -36247
$ws.Range("H96").Value = 30000
$ws.Range("I96").Value = 30000
$ws.Range("K96").Value = 30000
$ws.Range("M96").Value = -27254
$ws.Range("H134").Value = 1776773.4
$ws.Range("I134").Value = 1704201.2
$ws.Range("J134").Value = 2284778
$ws.Range("K134").Value = 5112603.6
$ws.Range("L134").Value = 6854334
$ws.Range("M134").Value = -5110068.6
$ws.Range("N134").Value = -6859404

$ws = $wb.Worksheets.Item("CRP")
# CRP: set 72 cell value(s)
$ws.Range("H7").Value = 243.76
$ws.Range("I7").Value = 126.9375
$ws.Range("J7").Value = 451.44446
$ws.Range("K7").Value = 126.9375
$ws.Range("L7").Value = 451.44446
$ws.Range("M7").Value = -13.9375
$ws.Range("N7").Value = -677.4444599999999
$ws.Range("H22").Value = 467.5
$ws.Range("J22").Value = 352.33334
$ws.Range("L22").Value = 352.33334
$ws.Range("N22").Value = -1052.33334
$ws.Range("H31").Value = 150126.7
$ws.Range("I31").Value = 215413.53
$ws.Range("J31").Value = 41315.332
$ws.Range("K31").Value = 215413.53
$ws.Range("L31").Value = 41315.332
$ws.Range("M31").Value = -215118.53
$ws.Range("N31").Value = -41905.332
$ws.Range("H34").Value = 150126.7
$ws.Range("I34").Value = 215413.53
$ws.Range("J34").Value = 41315.332
$ws.Range("K34").Value = 215413.53
$ws.Range("L34").Value = 41315.332
$ws.Range("M34").Value = -215211.53
$ws.Range("N34").Value = -41719.332
$ws.Range("H52").Value = 45000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H58").Value = 1766483.2
$ws.Range("I58").Value = 2058397.1
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 2058397.1
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -2058194.1
$ws.Range("N58").Value = -15406
$ws.Range("H86").Value = 8600.833000000001
$ws.Range("J86").Value = 8721
$ws.Range("L86").Value = 8721
$ws.Range("N86").Value = -10967
$ws.Range("H89").Value = 8600.833000000001
$ws.Range("J89").Value = 8721
$ws.Range("L89").Value = 43605
$ws.Range("N89").Value = -54837
$ws.Range("H105").Value = 31067.75
$ws.Range("I105").Value = 36331.4
$ws.Range("K105").Value = 36331.4
$ws.Range("M105").Value = -34584.4
$ws.Range("H132").Value = 33506642
$ws.Range("I132").Value = 38463504
$ws.Range("J132").Value = 1287030.8
$ws.Range("K132").Value = 115390512
$ws.Range("L132").Value = 3861092.4
$ws.Range("M132").Value = -115387982
$ws.Range("N132").Value = -3866152.4
$ws.Range("H134").Value = 6451.8335
$ws.Range("I134").Value = 6909.074
$ws.Range("K134").Value = 20727.222
$ws.Range("M134").Value = -18192.222
$ws.Range("H136").Value = 1766483.2
$ws.Range("I136").Value = 2058397.1
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 6175191.300000001
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -6172641.300000001
$ws.Range("N136").Value = -50100
$ws.Range("H139").Value = 89864.69
$ws.Range("J139").Value = 106428.43
$ws.Range("L139").Value = 106428.43
$ws.Range("N139").Value = -116708.43
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
# CRP: clear 2 cell(s)
$ws.Range("N52").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL: set 89 cell value(s)
$ws.Range("H2").Value = 93.181816
$ws.Range("I2").Value = 38.7
$ws.Range("K2").Value = 232.2
$ws.Range("M2").Value = -119.2
$ws.Range("H12").Value = 59.545456
$ws.Range("I12").Value = 41.333332
$ws.Range("J12").Value = 66.375
$ws.Range("K12").Value = 123.999996
$ws.Range("L12").Value = 199.125
$ws.Range("M12").Value = 49.000004
$ws.Range("N12").Value = -545.125
$ws.Range("H60").Value = 590.2222
$ws.Range("I60").Value = 110
$ws.Range("J60").Value = 1550.6666
$ws.Range("K60").Value = 330
$ws.Range("L60").Value = 4651.9998
$ws.Range("M60").Value = -79
$ws.Range("N60").Value = -5153.9998
$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("H89").Value = 200
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("H103").Value = 950.5
$ws.Range("I103").Value = 1200.8334
$ws.Range("K103").Value = 3602.5002
$ws.Range("M103").Value = -2723.5002
$ws.Range("H114").Value = 414
$ws.Range("I114").Value = 414
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 1242
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 2012
$ws.Range("H117").Value = 3599.077
$ws.Range("J117").Value = 5964.4287
$ws.Range("L117").Value = 17893.2861
$ws.Range("N117").Value = -24777.2861
$ws.Range("H118").Value = 4603.4
$ws.Range("J118").Value = 6447
$ws.Range("L118").Value = 19341
$ws.Range("N118").Value = -21827
$ws.Range("H122").Value = 15278590
$ws.Range("I122").Value = 766.55554
$ws.Range("J122").Value = 42778670
$ws.Range("K122").Value = 6898.99986
$ws.Range("L122").Value = 385008030
$ws.Range("M122").Value = -4448.99986
$ws.Range("N122").Value = -385012930
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H129").Value = 1733.5294
$ws.Range("J129").Value = 3999
$ws.Range("L129").Value = 11997
$ws.Range("N129").Value = -21997
$ws.Range("H136").Value = 17500
$ws.Range("I136").Value = 17500
$ws.Range("K136").Value = 52500
$ws.Range("M136").Value = -47400
$ws.Range("H137").Value = 1513.8572
$ws.Range("J137").Value = 1900
$ws.Range("L137").Value = 5700
$ws.Range("N137").Value = -15900
$ws.Range("H138").Value = 1390
$ws.Range("I138").Value = 1390
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4170
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 970
$ws.Range("H139").Value = 252433
$ws.Range("I139").Value = 252433
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 757299
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -752159
$ws.Range("H140").Value = 3210.8333
$ws.Range("I140").Value = 2003.75
$ws.Range("J140").Value = 5625
$ws.Range("K140").Value = 6011.25
$ws.Range("L140").Value = 16875
$ws.Range("M140").Value = -831.25
$ws.Range("N140").Value = -27235
$ws.Range("H141").Value = 4300
$ws.Range("I141").Value = 4300
$ws.Range("K141").Value = 12900
$ws.Range("M141").Value = -7720
# CUL: clear 7 cell(s)
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N114").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM: set 26 cell value(s)
$ws.Range("H52").Value = 34737.25
$ws.Range("J52").Value = 34737.25
$ws.Range("L52").Value = 34737.25
$ws.Range("N52").Value = -35255.25
$ws.Range("H80").Value = 247033.05
$ws.Range("I80").Value = 396840.94
$ws.Range("K80").Value = 396840.94
$ws.Range("M80").Value = -395842.94
$ws.Range("H83").Value = 247033.05
$ws.Range("I83").Value = 396840.94
$ws.Range("K83").Value = 1984204.7
$ws.Range("M83").Value = -1979212.7
$ws.Range("H97").Value = 1803.3914
$ws.Range("I97").Value = 1563.2106
$ws.Range("J97").Value = 2944.25
$ws.Range("K97").Value = 1563.2106
$ws.Range("L97").Value = 2944.25
$ws.Range("M97").Value = -1067.2106
$ws.Range("N97").Value = -3936.25
$ws.Range("H132").Value = 36150852
$ws.Range("I132").Value = 37489224
$ws.Range("J132").Value = 14799
$ws.Range("K132").Value = 112467672
$ws.Range("L132").Value = 44397
$ws.Range("M132").Value = -112465142
$ws.Range("N132").Value = -49457

$ws = $wb.Worksheets.Item("LTW")
# LTW: set 55 cell value(s)
$ws.Range("H7").Value = 4137.7856
$ws.Range("I7").Value = 3994.0833
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 3994.0833
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -3882.0833
$ws.Range("N7").Value = -5224
$ws.Range("H68").Value = 3610.8
$ws.Range("I68").Value = 4019.8
$ws.Range("K68").Value = 4019.8
$ws.Range("M68").Value = -3270.8
$ws.Range("H71").Value = 3610.8
$ws.Range("I71").Value = 4019.8
$ws.Range("K71").Value = 20099
$ws.Range("M71").Value = -16355
$ws.Range("H82").Value = 350
$ws.Range("I82").Value = 200
$ws.Range("J82").Value = 500
$ws.Range("K82").Value = 200
$ws.Range("L82").Value = 500
$ws.Range("M82").Value = 161
$ws.Range("N82").Value = -1222
$ws.Range("H85").Value = 350
$ws.Range("I85").Value = 200
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 200
$ws.Range("L85").Value = 500
$ws.Range("M85").Value = 1048
$ws.Range("N85").Value = -2996
$ws.Range("H100").Value = 27501
$ws.Range("H126").Value = 4137.7856
$ws.Range("I126").Value = 3994.0833
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 11982.2499
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -9512.249899999999
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 4352731
$ws.Range("I132").Value = 4973121.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 14919364.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -14916834.5
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 98225.234
$ws.Range("I136").Value = 2849.8333
$ws.Range("J136").Value = 179975.58
$ws.Range("K136").Value = 8549.499899999999
$ws.Range("L136").Value = 539926.74
$ws.Range("M136").Value = -5999.499899999999
$ws.Range("N136").Value = -545026.74
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("WVR")
# WVR: set 36 cell value(s)
$ws.Range("H81").Value = 6174.75
$ws.Range("I81").Value = 5974.75
$ws.Range("J81").Value = 6374.75
$ws.Range("K81").Value = 11949.5
$ws.Range("L81").Value = 12749.5
$ws.Range("M81").Value = -10888.5
$ws.Range("N81").Value = -14871.5
$ws.Range("H84").Value = 6174.75
$ws.Range("I84").Value = 5974.75
$ws.Range("J84").Value = 6374.75
$ws.Range("K84").Value = 59747.5
$ws.Range("L84").Value = 63747.5
$ws.Range("M84").Value = -54443.5
$ws.Range("N84").Value = -74355.5
$ws.Range("H123").Value = 47497.5
$ws.Range("J123").Value = 47497.5
$ws.Range("L123").Value = 47497.5
$ws.Range("N123").Value = -57297.5
$ws.Range("H126").Value = 6942.909
$ws.Range("I126").Value = 6995.5
$ws.Range("J126").Value = 6879.8
$ws.Range("K126").Value = 20986.5
$ws.Range("L126").Value = 20639.4
$ws.Range("M126").Value = -18516.5
$ws.Range("N126").Value = -25579.4
$ws.Range("H132").Value = 4794883.5
$ws.Range("I132").Value = 5441678.5
$ws.Range("J132").Value = 8599.799999999999
$ws.Range("K132").Value = 16325035.5
$ws.Range("L132").Value = 25799.4
$ws.Range("M132").Value = -16322505.5
$ws.Range("N132").Value = -30859.4
$ws.Range("H136").Value = 19246
$ws.Range("I136").Value = 22778
$ws.Range("K136").Value = 68334
$ws.Range("M136").Value = -65784
